# Added analysis of expert data on the "All Data" sheet.
# Shifts the second trial-of-three data block (I/K/L/M/N, originally in
# rows 6-10) up by one row into rows 6-8, then adds Average / Systematic
# Error / Random Error summary rows (10, 12, 13) mirroring the existing
# summary rows 34/36/37 at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: drop the old "Measurement" label cell (J6), shift I/K/L/M/N up ---
$ws.Range("J6").ClearContents()
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = 9.2
$ws.Range("L6").Value = 19.8
$ws.Range("M6").Value = 49.7
$ws.Range("N6").Value = 194.1

# --- Row 7 ---
$ws.Range("I7").Value = 2
$ws.Range("K7").Value = 9.8
$ws.Range("L7").Value = 20.5
$ws.Range("M7").Value = 49.1
$ws.Range("N7").Value = 198.8

# --- Row 8 ---
$ws.Range("I8").Value = 3
$ws.Range("K8").Value = 10
$ws.Range("L8").Value = 20.3
$ws.Range("M8").Value = 50.6
$ws.Range("N8").Value = 200.4

# --- Row 9: this data moved up into rows 6-8, so clear it here ---
$ws.Range("I9:N9").ClearContents()

# --- Row 10: becomes the "Average in µG" summary row ---
$ws.Range("I10").Value = "Average in µG"
$ws.Range("K10").Formula = "=AVERAGE(K3:K8)"
$ws.Range("L10").Formula = "=AVERAGE(L3:L8)"
$ws.Range("M10").Formula = "=AVERAGE(M3:M8)"
$ws.Range("N10").Formula = "=AVERAGE(N3:N8)"
$ws.Range("K10:N10").NumberFormat = "0.00\ "

# --- Row 12: "Systematic Error" summary row ---
$ws.Range("I12").Value = "Systematic Error"
$ws.Range("K12").Formula = "=((K10-10)/10)"
$ws.Range("L12").Formula = "=(L10-20)/20"
$ws.Range("M12").Formula = "=(M10-50)/50"
$ws.Range("N12").Formula = "=(N10-200)/200"
$ws.Range("K12:N12").NumberFormat = "0.00%"

# --- Row 13: "Random Error" summary row ---
$ws.Range("I13").Value = "Random Error"
$ws.Range("K13").Formula = "=(K5-K10)"
$ws.Range("L13").Formula = "=(L5-L10)"
$ws.Range("M13").Formula = "=(M8-M10)"
$ws.Range("N13").Formula = "=(N8-N10)"
$ws.Range("K13:N13").NumberFormat = "0.00\ "

# --- Move the selection cursor to where it ended up after the edit ---
$ws.Range("M7").Select()
